# Apply cryptos list price/volume updates (and two coin row swaps) per the
# Thu Jun  6 09:00:28 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.913.93"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "3.846.11"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "704.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "3.844.49"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "4.494.25"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "3.819.43"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "70.942.03"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.29%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.183"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").Value = "3.801.68"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +6.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("E42").Value = "  -5.72%  "
$ws.Range("E45").Value = "  -4.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.299"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "412.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.94%  "
